$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 8: fiscal-period labels -> drop oldest (1396/12), shift left, append newest (1401/12) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Header row 9: publish dates -> drop oldest, shift left, append newest ---
$ws.Range("D9").Value = "1399-03-13 (9)"
$ws.Range("E9").Value = "1400-03-11 (9)"
$ws.Range("F9").Value = "1401-03-11 (10)"
$ws.Range("G9").Value = "1402-02-27 (7)"

# H9's new text "1402-02-27" looks exactly like an ISO date with no other
# characters around it, so a plain .Value/.Formula assignment gets
# auto-converted to a date serial by Excel's input parser (and would also
# force a new NumberFormat style onto the cell). Route it through a text
# formula, then flatten the formula to a literal value via copy/paste
# special so the stored cell stays a plain string with the original style.
$ws.Range("H9").Formula = '="1402-02-27"'
$ws.Range("H9").Copy()
$ws.Range("H9").PasteSpecial(-4163)

# --- Data rows: shift each row's D:H values one column left, and put the new year's figure in H ---

# Row 11 - فروش (Sales)
$ws.Range("D11").Value = 11421
$ws.Range("E11").Value = 13539
$ws.Range("F11").Value = 10796
$ws.Range("G11").Value = 16287
$ws.Range("H11").Value = 27827

# Row 12 - بهای تمام شده کالای فروش رفته (COGS)
$ws.Range("D12").Value = -7850
$ws.Range("E12").Value = -8277
$ws.Range("F12").Value = -5309
$ws.Range("G12").Value = -8064
$ws.Range("H12").Value = -16737

# Row 13 - سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 3571
$ws.Range("E13").Value = 5262
$ws.Range("F13").Value = 5487
$ws.Range("G13").Value = 8223
$ws.Range("H13").Value = 11090

# Row 14 - هزینه های عمومی, اداری و تشکیلاتی
$ws.Range("D14").Value = -510
$ws.Range("E14").Value = -555
$ws.Range("F14").Value = -516
$ws.Range("G14").Value = -674
$ws.Range("H14").Value = -606

# Row 16 - خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = 98
$ws.Range("E16").Value = 192
$ws.Range("F16").Value = 69
$ws.Range("G16").Value = 153
$ws.Range("H16").Value = 81

# Row 17 - سود (زیان) عملیاتی
$ws.Range("D17").Value = 3158
$ws.Range("E17").Value = 4899
$ws.Range("F17").Value = 5040
$ws.Range("G17").Value = 7702
$ws.Range("H17").Value = 10564

# Row 18 - هزینه های مالی
$ws.Range("D18").Value = -986
$ws.Range("E18").Value = -1222
$ws.Range("F18").Value = -688
$ws.Range("G18").Value = -881
$ws.Range("H18").Value = -956

# Row 19 - خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = -55
$ws.Range("E19").Value = -948
$ws.Range("F19").Value = -310
$ws.Range("G19").Value = -127
$ws.Range("H19").Value = -7

# Row 20 - سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 2118
$ws.Range("E20").Value = 2729
$ws.Range("F20").Value = 4043
$ws.Range("G20").Value = 6694
$ws.Range("H20").Value = 9601

# Row 21 - مالیات
$ws.Range("D21").Value = -523
$ws.Range("E21").Value = -521
$ws.Range("F21").Value = -487
$ws.Range("G21").Value = -943
$ws.Range("H21").Value = -1422

# Row 22 - سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 1595
$ws.Range("E22").Value = 2208
$ws.Range("F22").Value = 3556
$ws.Range("G22").Value = 5751
$ws.Range("H22").Value = 8179

# Row 24 - سود (زیان) خالص
$ws.Range("D24").Value = 1595
$ws.Range("E24").Value = 2208
$ws.Range("F24").Value = 3556
$ws.Range("G24").Value = 5751
$ws.Range("H24").Value = 8179

# Row 26 - سرمایه
$ws.Range("D26").Value = 2521
$ws.Range("E26").Value = 1988
$ws.Range("F26").Value = 1128
$ws.Range("G26").Value = 1724
$ws.Range("H26").Value = 2550
